$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-24 Wednesday", "2024-07-25 Thursday"),
    @("140×6=840", "676×9=6084"),
    @("345×6=2070", "423×3=1269"),
    @("253×2=506", "816×3=2448"),
    @("394×2=788", "483×7=3381"),
    @("393×2=786", "133×8=1064"),
    @("895×7=6265", "720×2=1440"),
    @("522×4=2088", "521×6=3126"),
    @("690×3=2070", "400×2=800"),
    @("353×2=706", "115×3=345"),
    @("160×5=800", "571×9=5139"),
    @("645×4=2580", "640×4=2560"),
    @("906×7=6342", "812×5=4060"),
    @("531×2=1062", "221×9=1989"),
    @("487×4=1948", "586×6=3516"),
    @("726×5=3630", "797×2=1594"),
    @("318×5=1590", "316×3=948"),
    @("869×4=3476", "806×4=3224"),
    @("216×9=1944", "857×2=1714"),
    @("891×8=7128", "598×5=2990"),
    @("376×3=1128", "944×2=1888"),
    @("517×2=1034", "116×2=232"),
    @("358×2=716", "218×8=1744"),
    @("537×3=1611", "752×9=6768"),
    @("263×2=526", "939×7=6573"),
    @("690×7=4830", "666×3=1998")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
